$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.781.48'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.679.02'
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.46'
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.14'
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("E9").Value = '  -3.52%  '
$ws.Range("E10").Value = '  -3.47%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.42'
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.366'
$ws.Range("E12").Value = '  -3.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.149.78'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.19'
$ws.Range("E14").Value = '  -2.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.675.55'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("E16").Value = '  -2.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.677.06'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.77'
$ws.Range("E18").Value = '  -5.98%  '
$ws.Range("E19").Value = '  -3.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.73'
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.19'
$ws.Range("E21").Value = '  -5.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.509'
$ws.Range("E23").Value = '  -2.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.03'
$ws.Range("E24").Value = '  -2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.12'
$ws.Range("E27").Value = '  -2.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.39'
$ws.Range("E28").Value = '  +4.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0847'
$ws.Range("E29").Value = '  -6.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.93'
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.28'
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.85'
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.39'
$ws.Range("E36").Value = '  -3.17%  '
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '340.01'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.13'
$ws.Range("E39").Value = '  -2.48%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.927'
$ws.Range("E40").Value = '  -4.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.97'
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.34'
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.76'
$ws.Range("E43").Value = '  -5.12%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.615'
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.05'
$ws.Range("E45").Value = '  -4.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0552'
$ws.Range("E47").Value = '  -5.35%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.99'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0968'
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.86'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0240'
$ws.Range("E51").Value = '  -3.79%  '
